# Insert a new weekly price record for "Ajo" (Chino / Primera) dated
# 2022-02-04 (serial 44596) just above the current row 125, pushing the
# existing rows 125-191 down by one (to 126-192).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("125:125").Insert()

$ws.Range("A125").Value = 7
$ws.Range("B125").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C125").Value = "Ñuble"
$ws.Range("D125").Value = 44596
$ws.Range("E125").Value = 16
$ws.Range("F125").Value = 100112003
$ws.Range("G125").Value = "Ajo"
$ws.Range("H125").Value = "Chino"
$ws.Range("I125").Value = "Primera"
$ws.Range("J125").Value = 120
$ws.Range("K125").Value = 19000
$ws.Range("L125").Value = 20000
$ws.Range("M125").Value = 19500
$ws.Range("N125").Value = "$/caja 10 kilos"
$ws.Range("O125").Value = "China"
$ws.Range("P125").Value = 1950
$ws.Range("Q125").Value = 10
$ws.Range("R125").Value = "Hortaliza"
